# Update "想去人数" (interest count) figures in column F for the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets, matching the
# freshly generated data snapshot (gh-pages output @ 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 4668
$ws1.Range("F7").Value  = 583
$ws1.Range("F9").Value  = 274
$ws1.Range("F10").Value = 637
$ws1.Range("F15").Value = 1784
$ws1.Range("F18").Value = 1624
$ws1.Range("F19").Value = 13
$ws1.Range("F27").Value = 107
$ws1.Range("F30").Value = 32
$ws1.Range("F32").Value = 3886
$ws1.Range("F36").Value = 1020
$ws1.Range("F38").Value = 1858

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 4668
$ws4.Range("F7").Value  = 583
$ws4.Range("F9").Value  = 274
$ws4.Range("F10").Value = 637
$ws4.Range("F16").Value = 1784
$ws4.Range("F19").Value = 1624
$ws4.Range("F20").Value = 13
$ws4.Range("F28").Value = 107
$ws4.Range("F31").Value = 32
$ws4.Range("F33").Value = 3886
$ws4.Range("F38").Value = 1020
$ws4.Range("F40").Value = 1858
